$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")
$ws.Columns("P").Delete()
$ws.Columns("O").Delete()
$ws.Columns("N").Delete()
$ws.Columns("M").Delete()
$ws.Columns("F").Delete()
$ws.Columns("C").Delete()
for ($i=0; $i -lt 12; $i++) {
    $ws.Cells.Item(9, $i+2).Value2 = $i+1
}
$ws.Columns("H:M").AutoFit()
